$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project5186")

$ws.Range("B8").Value  = "T?;0=0;1=1;2=?"
$ws.Range("B12").Value = "T?;0=0;1=1;2=?"
$ws.Range("B19").Value = "T?;0=0;1=1;2=2;3=?"
$ws.Range("B21").Value = "T?;0=0;1=1;2=?"
$ws.Range("B23").Value = "N;0=0;1=1;2=?"
$ws.Range("B26").Value = "T?;0=0;1=1;2=?"
$ws.Range("B31").Value = "T?;0=0;1=1;2=?"
$ws.Range("B34").Value = "T?;0=0;1=1;2=?"
$ws.Range("B36").Value = "N;0=0;1=1;2=?"
$ws.Range("B38").Value = "N;0=0;1=1;2=?"
$ws.Range("B41").Value = "T?;0=0;1=1;2=?"

$ws.Range("B38").Select() | Out-Null
